$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 0.9999970461623656
$block1[0,2] = 1.030263256528299
$block1[0,3] = 1.003512603051345
$block1[0,4] = 0.9980539377174651
$block1[1,0] = 1.02
$block1[1,1] = 1.002190696077171
$block1[1,2] = 1.030606612076188
$block1[1,3] = 1.005418250165748
$block1[1,4] = 1.000925108801209
$block1[2,0] = 1.02
$block1[2,1] = 1.003602290855779
$block1[2,2] = 1.030828240658591
$block1[2,3] = 1.006644382404554
$block1[2,4] = 1.002773591153557
$block1[3,0] = 1.02
$block1[3,1] = 1.004193894425286
$block1[3,2] = 1.030921280219693
$block1[3,3] = 1.007158223888107
$block1[3,4] = 1.003548520839981
$block1[4,0] = 1.02
$block1[4,1] = 1.004293121253613
$block1[4,2] = 1.030936894093018
$block1[4,3] = 1.007244405997112
$block1[4,4] = 1.003678509513614
$block1[5,0] = 1.02
$block1[5,1] = 1.003610203023348
$block1[5,2] = 1.030829484383891
$block1[5,3] = 1.006651254710766
$block1[5,4] = 1.002783954247485
$block1[6,0] = 1.02
$block1[6,1] = 1.000740052902658
$block1[6,2] = 1.030379406564799
$block1[6,3] = 1.004158087951371
$block1[6,4] = 0.9990262421512628
$block1[7,0] = 1.02
$block1[7,1] = 0.9956202764081897
$block1[7,2] = 1.029582245551908
$block1[7,3] = 0.9997098096606482
$block1[7,4] = 0.9923299142558761
$block1[8,0] = 1.02
$block1[8,1] = 0.9921622425845256
$block1[8,2] = 1.029048220613477
$block1[8,3] = 0.9967047873794241
$block1[8,4] = 0.9878110444644356
$block1[9,0] = 1.02
$block1[9,1] = 0.9906535521406715
$block1[9,2] = 1.028816400681361
$block1[9,3] = 0.9953936374089645
$block1[9,4] = 0.9858403664321082
$block1[10,0] = 1.02
$block1[10,1] = 0.9900913939142384
$block1[10,2] = 1.028730207336771
$block1[10,3] = 0.9949050722402973
$block1[10,4] = 0.9851061820500319
$block1[11,0] = 1.02
$block1[11,1] = 0.9902120595381855
$block1[11,2] = 1.028748699911505
$block1[11,3] = 0.9950099418895325
$block1[11,4] = 0.9852637674690576
$block1[12,0] = 1.02
$block1[12,1] = 0.9906071202911284
$block1[12,2] = 1.028809277633772
$block1[12,3] = 0.9953532842849065
$block1[12,4] = 0.9857797236426981
$block1[13,0] = 1.02
$block1[13,1] = 0.9908502947745117
$block1[13,2] = 1.028846590368057
$block1[13,3] = 0.9955646225320298
$block1[13,4] = 0.9860973288618624
$block1[14,0] = 1.02
$block1[14,1] = 0.9922621254230287
$block1[14,2] = 1.029063593663636
$block1[14,3] = 0.9967915900936591
$block1[14,4] = 0.9879415300188076
$block1[15,0] = 1.02
$block1[15,1] = 0.9931446518951077
$block1[15,2] = 1.029199559676678
$block1[15,3] = 0.9975585341355973
$block1[15,4] = 0.989094545834366
$block1[16,0] = 1.02
$block1[16,1] = 0.993658324721283
$block1[16,2] = 1.029278809869972
$block1[16,3] = 0.9980049224938602
$block1[16,4] = 0.9897657377869169
$block1[17,0] = 1.02
$block1[17,1] = 0.9938332907999199
$block1[17,2] = 1.029305822478735
$block1[17,3] = 0.9981569685386646
$block1[17,4] = 0.9899943717553462
$block1[18,0] = 1.02
$block1[18,1] = 0.9930500783226616
$block1[18,2] = 1.029184977633817
$block1[18,3] = 0.9974763476563392
$block1[18,4] = 0.9889709776493107
$block1[19,0] = 1.02
$block1[19,1] = 0.9904908339126361
$block1[19,2] = 1.028791441332969
$block1[19,3] = 0.9952522215828185
$block1[19,4] = 0.9856278484838727
$block1[20,0] = 1.02
$block1[20,1] = 0.9888715021031147
$block1[20,2] = 1.0285435183423
$block1[20,3] = 0.9938448567504581
$block1[20,4] = 0.9835131937094171
$block1[21,0] = 1.02
$block1[21,1] = 0.9897309308844535
$block1[21,2] = 1.028674992704237
$block1[21,3] = 0.9945917945827535
$block1[21,4] = 0.9846354452697469
$block1[22,0] = 1.02
$block1[22,1] = 0.9930928154140612
$block1[22,2] = 1.029191566806253
$block1[22,3] = 0.9975134871437007
$block1[22,4] = 0.9890268169374212
$block1[23,0] = 1.02
$block1[23,1] = 0.9969515404821367
$block1[23,2] = 1.029788797171083
$block1[23,3] = 1.000866570495859
$block1[23,4] = 0.9940703804778075
$ws.Range("B2:F25").Value = $block1

$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 1.02955528997503
$block2[0,1] = 1.005312350954267
$block2[0,2] = 1.033074774661431
$block2[0,3] = 1.006403886417997
$block2[0,4] = 1.000962071010937
$block2[0,5] = 1.005586813710726
$block2[1,0] = 1.029546965970127
$block2[1,1] = 1.007127342732227
$block2[1,2] = 1.033227698754684
$block2[1,3] = 1.008108944823934
$block2[1,4] = 1.003628619941815
$block2[1,5] = 1.006228286426888
$block2[2,0] = 1.029539753632458
$block2[2,1] = 1.008294031562883
$block2[2,2] = 1.0333253549493
$block2[2,3] = 1.009204923891012
$block2[2,4] = 1.005344591826161
$block2[2,5] = 1.006639578476929
$block2[3,0] = 1.029536283262567
$block2[3,1] = 1.008782696506916
$block2[3,2] = 1.0333660975551
$block2[3,3] = 1.009663960881819
$block2[3,4] = 1.006063784403245
$block2[3,5] = 1.006811593515251
$block2[4,0] = 1.029535674837201
$block2[4,1] = 1.008864640454659
$block2[4,2] = 1.033372920065691
$block2[4,3] = 1.009740935846062
$block2[4,4] = 1.006184412705869
$block2[4,5] = 1.006840423675391
$block2[5,0] = 1.029539708984831
$block2[5,1] = 1.008300568190578
$block2[5,2] = 1.033325900582437
$block2[5,3] = 1.009211064246002
$block2[5,4] = 1.005354210280489
$block2[5,5] = 1.006641880437511
$block2[6,0] = 1.02955285428455
$block2[6,1] = 1.005927361266509
$block2[6,2] = 1.033126723317438
$block2[6,3] = 1.006981654404119
$block2[6,4] = 1.001865241835927
$block2[6,5] = 1.005804394172819
$block2[7,0] = 1.02956208100376
$block2[7,1] = 1.001684422044947
$block2[7,2] = 1.032765901988229
$block2[7,3] = 1.002995500289088
$block2[7,4] = 0.9956418554368396
$block2[7,5] = 1.004299033365596
$block2[8,0] = 1.029558927936093
$block2[8,1] = 0.9988121622532348
$block2[8,2] = 1.032518842347417
$block2[8,3] = 1.000296919877218
$block2[8,4] = 0.9914381395819487
$block2[8,5] = 1.003274668619613
$block2[9,0] = 1.029555369739793
$block2[9,1] = 0.9975575053047278
$block2[9,2] = 1.032410340936441
$block2[9,3] = 0.9991181050659224
$block2[9,4] = 0.9896039559150833
$block2[9,5] = 1.002825974291096
$block2[10,0] = 1.029553719719501
$block2[10,1] = 0.9970897735031714
$block2[10,2] = 1.032369811708945
$block2[10,3] = 0.9986786442991347
$block2[10,4] = 0.9889204816658863
$block2[10,5] = 1.002658519692851
$block2[11,0] = 1.029554088495152
$block2[11,1] = 0.9971901811635615
$block2[11,2] = 1.032378515607504
$block2[11,3] = 0.9987729831400927
$block2[11,4] = 0.9890671890455449
$block2[11,5] = 1.00269447529006
$block2[12,0] = 1.029555240040406
$block2[12,1] = 0.9975188773544146
$block2[12,2] = 1.032406995397466
$block2[12,3] = 0.9990818119466078
$block2[12,4] = 0.9895475046040829
$block2[12,5] = 1.002812148652134
$block2[13,0] = 1.029555906069242
$block2[13,1] = 0.9977211715879997
$block2[13,2] = 1.032424512714591
$block2[13,3] = 0.9992718785883367
$block2[13,4] = 0.9898431519607224
$block2[13,5] = 1.002884545935586
$block2[14,0] = 1.029559117975006
$block2[14,1] = 0.9988951946573905
$block2[14,2] = 1.032526011298701
$block2[14,3] = 1.00037493267123
$block2[14,4] = 0.9915595675619321
$block2[14,5] = 1.003304337241511
$block2[15,0] = 1.029560546404202
$block2[15,1] = 0.9996286609749276
$block2[15,2] = 1.032589272239751
$block2[15,3] = 1.001064055618971
$block2[15,4] = 0.992632438877268
$block2[15,5] = 1.003566273402792
$block2[16,0] = 1.029561167926103
$block2[16,1] = 1.000055426526088
$block2[16,2] = 1.032626024211051
$block2[16,3] = 1.001465017693937
$block2[16,4] = 0.9932568864796915
$block2[16,5] = 1.003718561842482
$block2[17,0] = 1.029561343914426
$block2[17,1] = 1.000200765487866
$block2[17,2] = 1.032638530682611
$block2[17,3] = 1.001601568606567
$block2[17,4] = 0.9934695821076397
$block2[17,5] = 1.003770404919007
$block2[18,0] = 1.029560415031907
$block2[18,1] = 0.9995500762588492
$block2[18,2] = 1.032582500137903
$block2[18,3] = 1.00099022216119
$block2[18,4] = 0.9925174691326852
$block2[18,5] = 1.003538221405817
$block2[19,0] = 1.029554909994161
$block2[19,1] = 0.9974221317641361
$block2[19,2] = 1.032398615060048
$block2[19,3] = 0.9989909140153368
$block2[19,4] = 0.9894061243963045
$block2[19,5] = 1.002777518719308
$block2[20,0] = 1.029549549582982
$block2[20,1] = 0.9960743677744542
$block2[20,2] = 1.032281686918612
$block2[20,3] = 0.9977246088080334
$block2[20,4] = 0.9874372626241351
$block2[20,5] = 1.002294657611214
$block2[21,0] = 1.029552570889732
$block2[21,1] = 0.9967897932096135
$block2[21,2] = 1.032343796520178
$block2[21,3] = 0.998396795041804
$block2[21,4] = 0.9884822189254497
$block2[21,5] = 1.002551071293827
$block2[22,0] = 1.029560475047559
$block2[22,1] = 0.9995855885604165
$block2[22,2] = 1.032585560613934
$block2[22,3] = 1.001023587384582
$block2[22,4] = 0.9925694231479877
$block2[22,5] = 1.003550898422906
$block2[23,0] = 1.029561340424493
$block2[23,1] = 1.002788820867278
$block2[23,2] = 1.032860338570589
$block2[23,3] = 1.004033090616584
$block2[23,4] = 0.9972601003021522
$block2[23,5] = 1.004691803431707
$ws.Range("I2:N25").Value = $block2

Write-Host "Applied vm_pu.xlsx updates"